$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to retain their existing text type
# (these are inline-string cells holding numeric-looking / percentage text,
# not real numbers) while writing the new values, then restore the default
# "Normal" style so no stray style index is left on the cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.957.44"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.635.09"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "214.22"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("D10").Value = "18.52"
$ws.Range("E10").Value = "  -5.77%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "1.862.60"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "1.634.97"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D16").Value = "25.963.41"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "0.0₃0742"
$ws.Range("E17").Value = "  -3.05%  "
$ws.Range("D18").Value = "61.72"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "190.70"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").Value = "9.69"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "1.78"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "143.41"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("D32").Value = "3.16"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.49"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.41"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("D36").Value = "1.135.85"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D38").Value = "2.44"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "0.521"
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").Value = "98.35"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").Value = "1.772.56"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  -5.00%  "
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "55.20"
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D50").Value = "7.51"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("E51").Value = "  -0.01%  "

$ws.Range("D2:E51").Style = "Normal"
